$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'68.437.35"
$ws.Cells.Item(2, 5).Value = "'  -1.08%  "
$ws.Cells.Item(3, 4).Value = "'2.431.39"
$ws.Cells.Item(3, 5).Value = "'  -1.75%  "
$ws.Cells.Item(4, 5).Value = "'  -0.02%  "
$ws.Cells.Item(5, 4).Value = "'557.90"
$ws.Cells.Item(5, 5).Value = "'  -0.66%  "
$ws.Cells.Item(6, 4).Value = "'160.25"
$ws.Cells.Item(6, 5).Value = "'  -1.77%  "
$ws.Cells.Item(7, 5).Value = "'  -0.01%  "
$ws.Cells.Item(8, 5).Value = "'  +0.03%  "
$ws.Cells.Item(9, 5).Value = "'  +6.72%  "
$ws.Cells.Item(10, 5).Value = "'  -1.78%  "
$ws.Cells.Item(11, 5).Value = "'  -0.71%  "
$ws.Cells.Item(12, 4).Value = "'4.64"
$ws.Cells.Item(12, 5).Value = "'  -5.06%  "
$ws.Cells.Item(13, 4).Value = "'68.295.10"
$ws.Cells.Item(13, 5).Value = "'  -1.11%  "
$ws.Cells.Item(14, 4).Value = "'2.871.95"
$ws.Cells.Item(14, 5).Value = "'  -1.16%  "
$ws.Cells.Item(15, 5).Value = "'  +1.82%  "
$ws.Cells.Item(16, 4).Value = "'23.07"
$ws.Cells.Item(16, 5).Value = "'  -2.89%  "
$ws.Cells.Item(17, 4).Value = "'2.425.85"
$ws.Cells.Item(17, 5).Value = "'  -2.05%  "
$ws.Cells.Item(18, 4).Value = "'10.46"
$ws.Cells.Item(18, 5).Value = "'  -3.10%  "
$ws.Cells.Item(19, 4).Value = "'335.11"
$ws.Cells.Item(19, 5).Value = "'  -1.28%  "
$ws.Cells.Item(20, 4).Value = "'6.90"
$ws.Cells.Item(20, 5).Value = "'  -1.69%  "
$ws.Cells.Item(21, 5).Value = "'  +0.39%  "
$ws.Cells.Item(22, 4).Value = "'1.91"
$ws.Cells.Item(22, 5).Value = "'  +0.61%  "
$ws.Cells.Item(23, 5).Value = "'  +0.00%  "
$ws.Cells.Item(24, 4).Value = "'66.63"
$ws.Cells.Item(24, 5).Value = "'  -1.04%  "
$ws.Cells.Item(25, 4).Value = "'3.68"
$ws.Cells.Item(25, 5).Value = "'  -0.47%  "
$ws.Cells.Item(26, 4).Value = "'2.551.02"
$ws.Cells.Item(26, 5).Value = "'  -1.93%  "
$ws.Cells.Item(27, 4).Value = "'8.21"
$ws.Cells.Item(27, 5).Value = "'  -1.19%  "
$ws.Cells.Item(28, 4).Value = "'0.0₃0817"
$ws.Cells.Item(28, 5).Value = "'  -1.01%  "
$ws.Cells.Item(29, 4).Value = "'7.15"
$ws.Cells.Item(29, 5).Value = "'  -0.99%  "
$ws.Cells.Item(30, 5).Value = "'  +0.04%  "
$ws.Cells.Item(31, 4).Value = "'426.06"
$ws.Cells.Item(31, 5).Value = "'  -1.55%  "
$ws.Cells.Item(32, 4).Value = "'1.15"
$ws.Cells.Item(32, 5).Value = "'  -0.26%  "
$ws.Cells.Item(33, 5).Value = "'  -1.40%  "
$ws.Cells.Item(34, 4).Value = "'158.79"
$ws.Cells.Item(34, 5).Value = "'  +0.69%  "
$ws.Cells.Item(35, 5).Value = "'  -0.14%  "
$ws.Cells.Item(36, 5).Value = "'  -0.04%  "
$ws.Cells.Item(37, 5).Value = "'  +0.30%  "
$ws.Cells.Item(38, 5).Value = "'  -4.66%  "
$ws.Cells.Item(39, 5).Value = "'  -1.72%  "
$ws.Cells.Item(40, 4).Value = "'4.33"
$ws.Cells.Item(40, 5).Value = "'  -2.88%  "
$ws.Cells.Item(41, 4).Value = "'1.49"
$ws.Cells.Item(41, 5).Value = "'  +0.51%  "
$ws.Cells.Item(42, 4).Value = "'1.08"
$ws.Cells.Item(42, 5).Value = "'  -1.07%  "
$ws.Cells.Item(43, 4).Value = "'133.19"
$ws.Cells.Item(43, 5).Value = "'  +0.04%  "
$ws.Cells.Item(44, 4).Value = "'2.02"
$ws.Cells.Item(44, 5).Value = "'  -2.87%  "
$ws.Cells.Item(45, 4).Value = "'3.34"
$ws.Cells.Item(45, 5).Value = "'  -0.64%  "
$ws.Cells.Item(46, 5).Value = "'  -0.48%  "
$ws.Cells.Item(47, 4).Value = "'0.481"
$ws.Cells.Item(47, 5).Value = "'  -1.08%  "
$ws.Cells.Item(48, 5).Value = "'  -1.40%  "
$ws.Cells.Item(49, 4).Value = "'0.0915"
$ws.Cells.Item(49, 5).Value = "'  -0.41%  "
$ws.Cells.Item(50, 5).Value = "'  -0.10%  "
$ws.Cells.Item(51, 5).Value = "'  -1.77%  "
